$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1464.579
$ws.Range("J17").Value = 1464.579
$ws.Range("L17").Value = 4393.737
$ws.Range("N17").Value = -4729.737
$ws.Range("H40").Value = 2089.0908
$ws.Range("J40").Value = 2796
$ws.Range("L40").Value = 2796
$ws.Range("N40").Value = -3146
$ws.Range("H57").Value = 24500
$ws.Range("J57").Value = 24500
$ws.Range("L57").Value = 73500
$ws.Range("N57").Value = -74498
$ws.Range("H64").Value = 4383.4
$ws.Range("I64").Value = 3349.8333
$ws.Range("J64").Value = 4641.7915
$ws.Range("K64").Value = 3349.8333
$ws.Range("L64").Value = 4641.7915
$ws.Range("M64").Value = -3101.8333
$ws.Range("N64").Value = -5137.7915
$ws.Range("H67").Value = 4383.4
$ws.Range("I67").Value = 3349.8333
$ws.Range("J67").Value = 4641.7915
$ws.Range("K67").Value = 3349.8333
$ws.Range("L67").Value = 4641.7915
$ws.Range("M67").Value = -2491.8333
$ws.Range("N67").Value = -6357.7915
$ws.Range("H74").Value = 4200
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("H77").Value = 4200
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("H98").Value = 7278.76
$ws.Range("I98").Value = 5274
$ws.Range("J98").Value = 15297.8
$ws.Range("K98").Value = 5274
$ws.Range("L98").Value = 15297.8
$ws.Range("M98").Value = -3776
$ws.Range("N98").Value = -18293.8
$ws.Range("H100").Value = 2274
$ws.Range("I100").Value = 1760
$ws.Range("J100").Value = 2788
$ws.Range("K100").Value = 1760
$ws.Range("L100").Value = 2788
$ws.Range("M100").Value = -1219
$ws.Range("N100").Value = -3870
$ws.Range("H112").Value = 1969.4375
$ws.Range("J112").Value = 2125.75
$ws.Range("L112").Value = 6377.25
$ws.Range("N112").Value = -8593.25
$ws.Range("H113").Value = 3888.75
$ws.Range("I113").Value = 3902.5
$ws.Range("J113").Value = 3884.1667
$ws.Range("K113").Value = 3902.5
$ws.Range("L113").Value = 3884.1667
$ws.Range("M113").Value = -648.5
$ws.Range("N113").Value = -10392.1667
$ws.Range("H122").Value = 7278.76
$ws.Range("I122").Value = 5274
$ws.Range("J122").Value = 15297.8
$ws.Range("K122").Value = 15822
$ws.Range("L122").Value = 45893.39999999999
$ws.Range("M122").Value = -13372
$ws.Range("N122").Value = -50793.39999999999
$ws.Range("H129").Value = 984.4318
$ws.Range("J129").Value = 1160.6061
$ws.Range("L129").Value = 3481.8183
$ws.Range("N129").Value = -13481.8183
$ws.Range("H132").Value = 1216.2954
$ws.Range("I132").Value = 1117.317
$ws.Range("K132").Value = 3351.951
$ws.Range("M132").Value = -821.951
$ws.Range("H138").Value = 2181.94
$ws.Range("I138").Value = 1472.5853
$ws.Range("J138").Value = 2674.8813
$ws.Range("K138").Value = 4417.7559
$ws.Range("L138").Value = 8024.6439
$ws.Range("M138").Value = 722.2440999999999
$ws.Range("N138").Value = -18304.6439
$ws.Range("M74").ClearContents()
$ws.Range("M77").ClearContents()

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H35").Value = 50000
$ws.Range("I35").Value = 50000
$ws.Range("K35").Value = 50000
$ws.Range("M35").Value = -49594
$ws.Range("H74").Value = 1688.4166
$ws.Range("I74").Value = 1440.5
$ws.Range("J74").Value = 2184.25
$ws.Range("K74").Value = 1440.5
$ws.Range("L74").Value = 2184.25
$ws.Range("M74").Value = -566.5
$ws.Range("N74").Value = -3932.25
$ws.Range("H77").Value = 1688.4166
$ws.Range("I77").Value = 1440.5
$ws.Range("J77").Value = 2184.25
$ws.Range("K77").Value = 7202.5
$ws.Range("L77").Value = 10921.25
$ws.Range("M77").Value = -2834.5
$ws.Range("N77").Value = -19657.25
$ws.Range("M97").Value = -157.75
$ws.Range("H122").Value = 2550.7144
$ws.Range("I122").Value = 2550.7144
$ws.Range("K122").Value = 7652.1432
$ws.Range("M122").Value = -5202.1432
$ws.Range("N97").ClearContents()

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 83778.57000000001
$ws.Range("I82").Value = 89527.836
$ws.Range("J82").Value = 49283
$ws.Range("K82").Value = 89527.836
$ws.Range("L82").Value = 49283
$ws.Range("M82").Value = -89144.836
$ws.Range("N82").Value = -50049
$ws.Range("H85").Value = 83778.57000000001
$ws.Range("I85").Value = 89527.836
$ws.Range("J85").Value = 49283
$ws.Range("K85").Value = 89527.836
$ws.Range("L85").Value = 49283
$ws.Range("M85").Value = -88201.836
$ws.Range("N85").Value = -51935
$ws.Range("H86").Value = 51906.4
$ws.Range("I86").Value = 1894
$ws.Range("J86").Value = 144786.58
$ws.Range("K86").Value = 1894
$ws.Range("L86").Value = 144786.58
$ws.Range("M86").Value = -771
$ws.Range("N86").Value = -147032.58
$ws.Range("H89").Value = 51906.4
$ws.Range("I89").Value = 1894
$ws.Range("J89").Value = 144786.58
$ws.Range("K89").Value = 9470
$ws.Range("L89").Value = 723932.8999999999
$ws.Range("M89").Value = -3854
$ws.Range("N89").Value = -735164.8999999999
$ws.Range("H102").Value = 142414.5
$ws.Range("I102").Value = 142414.5
$ws.Range("K102").Value = 142414.5
$ws.Range("M102").Value = -139169.5
$ws.Range("H119").Value = 24999.5
$ws.Range("J119").Value = 24999.5
$ws.Range("L119").Value = 24999.5
$ws.Range("N119").Value = -34675.5

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2513.551
$ws.Range("I31").Value = 1817.421
$ws.Range("K31").Value = 1817.421
$ws.Range("M31").Value = -1522.421
$ws.Range("H34").Value = 2513.551
$ws.Range("I34").Value = 1817.421
$ws.Range("K34").Value = 1817.421
$ws.Range("M34").Value = -1615.421
$ws.Range("H132").Value = 308780.8
$ws.Range("I132").Value = 423709.8
$ws.Range("J132").Value = 2303.5
$ws.Range("K132").Value = 1271129.4
$ws.Range("L132").Value = 6910.5
$ws.Range("M132").Value = -1268599.4
$ws.Range("N132").Value = -11970.5
$ws.Range("H134").Value = 1597.6296
$ws.Range("I134").Value = 1451.6364
$ws.Range("J134").Value = 2240
$ws.Range("K134").Value = 4354.9092
$ws.Range("L134").Value = 6720
$ws.Range("M134").Value = -1819.9092
$ws.Range("N134").Value = -11790

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1249.75
$ws.Range("I68").Value = 2000
$ws.Range("J68").Value = 999.6667
$ws.Range("K68").Value = 6000
$ws.Range("L68").Value = 2999.0001
$ws.Range("M68").Value = -5189
$ws.Range("N68").Value = -4621.0001
$ws.Range("H71").Value = 1249.75
$ws.Range("I71").Value = 2000
$ws.Range("J71").Value = 999.6667
$ws.Range("K71").Value = 18000
$ws.Range("L71").Value = 8997.0003
$ws.Range("M71").Value = -13944
$ws.Range("N71").Value = -17109.0003
$ws.Range("H131").Value = 14495785
$ws.Range("J131").Value = 16130753
$ws.Range("L131").Value = 48392259
$ws.Range("N131").Value = -48402339

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4363.7144
$ws.Range("I122").Value = 4007.36
$ws.Range("K122").Value = 12022.08
$ws.Range("M122").Value = -9572.08
$ws.Range("H123").Value = 33300.92
$ws.Range("J123").Value = 33300.92
$ws.Range("L123").Value = 33300.92
$ws.Range("N123").Value = -38200.92
$ws.Range("H140").Value = 44136.668
$ws.Range("J140").Value = 44136.668
$ws.Range("L140").Value = 44136.668
$ws.Range("N140").Value = -54496.668
$ws.Range("H141").Value = 38163.332
$ws.Range("J141").Value = 38163.332
$ws.Range("L141").Value = 38163.332
$ws.Range("N141").Value = -48523.332

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5611.9033
$ws.Range("I132").Value = 5826.577
$ws.Range("J132").Value = 4495.6
$ws.Range("K132").Value = 17479.731
$ws.Range("L132").Value = 13486.8
$ws.Range("M132").Value = -14949.731
$ws.Range("N132").Value = -18546.8
$ws.Range("H139").Value = 37581.125
$ws.Range("J139").Value = 37581.125
$ws.Range("L139").Value = 37581.125
$ws.Range("N139").Value = -47861.125

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 49450
$ws.Range("J68").Value = 49450
$ws.Range("L68").Value = 49450
$ws.Range("N68").Value = -51072
$ws.Range("H71").Value = 49450
$ws.Range("J71").Value = 49450
$ws.Range("L71").Value = 148350
$ws.Range("N71").Value = -156462
$ws.Range("H80").Value = 45266.668
$ws.Range("J80").Value = 45266.668
$ws.Range("L80").Value = 45266.668
$ws.Range("N80").Value = -47262.668
$ws.Range("H83").Value = 45266.668
$ws.Range("J83").Value = 45266.668
$ws.Range("L83").Value = 135800.004
$ws.Range("N83").Value = -145784.004
$ws.Range("H96").Value = 2000
$ws.Range("J96").Value = 2000
$ws.Range("L96").Value = 2000
$ws.Range("N96").Value = -4746
$ws.Range("H122").Value = 1099.1428
$ws.Range("I122").Value = 1138.8
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 3416.4
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -966.3999999999996
$ws.Range("N122").Value = -7900
